$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel
# auto-converting numeric-looking strings like "588.69" or "1.00" into
# floating point numbers, which would lose formatting / precision).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.943.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.492.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.42"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.145"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.98"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.943.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.48"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.843.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.482.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.91"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.592.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0900"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "500.28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.77"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.96"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.64"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.82"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0257"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.79%  "
